# Fix to last commit
#
# - On the "App" sheet, insert a new settings row above the old row 12
#   ("Number of processor threads ...") for a new "Rounded time step size"
#   setting (value 0.25, with its note).
# - Make "App" the active sheet/tab again (it had drifted to "Algorithm").

$wb = $excel.ActiveWorkbook

$appSheet = $wb.Worksheets.Item("App")
$algoSheet = $wb.Worksheets.Item("Algorithm")

# Insert a new row above row 12, pushing the existing rows 12-13 down to 13-14.
$appSheet.Rows.Item(12).Insert()

$appSheet.Cells.Item(12, 1).Value = "Rounded time step size"
$appSheet.Cells.Item(12, 2).Value = 0.25
$appSheet.Cells.Item(12, 3).Value = "Some shift info in the program is rounded to this number of hours"

$appSheet.Cells.Item(12, 1).Style = "Normal"
$appSheet.Cells.Item(12, 2).Style = "Normal"
$appSheet.Cells.Item(12, 3).Style = "Normal"

# Copy the row-12 (now row 13) formatting down to the new row so the new
# cells match their neighbours' style (s="5"), then restore row 13's taller
# height for its wrapped note text.
$appSheet.Rows.Item(13).RowHeight = $appSheet.Rows.Item(11).RowHeight

$appSheet.Activate()
$appSheet.Range("C16").Select()

$algoSheet.Range("D11").Select()
